# Applies the cryptos.xlsx price/volume(1h) refresh + row-27/28 swap + row-51 coin replacement
# described by the GitHub Actions "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Sheet, $Ref, $Text) {
    $range = $Sheet.Range($Ref)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.ClearFormats()
}

Set-TextCell $ws "D2" "51.747.80"
Set-TextCell $ws "E2" "  +1.30%  "
Set-TextCell $ws "D3" "3.065.06"
Set-TextCell $ws "E3" "  +3.67%  "
Set-TextCell $ws "D4" "1.00"
Set-TextCell $ws "E4" "  +0.10%  "
Set-TextCell $ws "D5" "384.41"
Set-TextCell $ws "E5" "  +1.20%  "
Set-TextCell $ws "D6" "103.45"
Set-TextCell $ws "E6" "  +1.38%  "
Set-TextCell $ws "E7" "  +0.41%  "
Set-TextCell $ws "E8" "  +0.00%  "
Set-TextCell $ws "D9" "0.589"
Set-TextCell $ws "E9" "  +0.05%  "
Set-TextCell $ws "D10" "37.20"
Set-TextCell $ws "E10" "  +2.33%  "
Set-TextCell $ws "E11" "  +0.43%  "
Set-TextCell $ws "E12" "  +1.22%  "
Set-TextCell $ws "D13" "3.534.30"
Set-TextCell $ws "E13" "  +3.43%  "
Set-TextCell $ws "D14" "18.74"
Set-TextCell $ws "E14" "  +2.53%  "
Set-TextCell $ws "D15" "7.78"
Set-TextCell $ws "E15" "  -0.05%  "
Set-TextCell $ws "D16" "3.076.45"
Set-TextCell $ws "E16" "  +4.25%  "
Set-TextCell $ws "D17" "0.983"
Set-TextCell $ws "E17" "  -1.09%  "
Set-TextCell $ws "D18" "10.54"
Set-TextCell $ws "E18" "  -5.36%  "
Set-TextCell $ws "D19" "51.802.64"
Set-TextCell $ws "E19" "  +1.29%  "
Set-TextCell $ws "D20" "3.15"
Set-TextCell $ws "E20" "  +0.51%  "
Set-TextCell $ws "D21" "12.52"
Set-TextCell $ws "E21" "  +1.57%  "
Set-TextCell $ws "D22" "0.0₃0965"
Set-TextCell $ws "E22" "  +0.60%  "
Set-TextCell $ws "D23" "70.17"
Set-TextCell $ws "E23" "  -0.31%  "
Set-TextCell $ws "D24" "268.80"
Set-TextCell $ws "E24" "  +0.69%  "
Set-TextCell $ws "D25" "3.17"
Set-TextCell $ws "E25" "  -3.30%  "
Set-TextCell $ws "D26" "8.54"
Set-TextCell $ws "E26" "  +9.20%  "
Set-TextCell $ws "B27" "Kaspa"
Set-TextCell $ws "C27" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws "D27" "0.173"
Set-TextCell $ws "E27" "  +5.25%  "
Set-TextCell $ws "B28" "EthereumClassic"
Set-TextCell $ws "C28" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws "D28" "27.07"
Set-TextCell $ws "E28" "  +4.84%  "
Set-TextCell $ws "E29" "  +2.26%  "
Set-TextCell $ws "E31" "  -2.42%  "
Set-TextCell $ws "E32" "  +0.65%  "
Set-TextCell $ws "D33" "34.62"
Set-TextCell $ws "E33" "  +1.12%  "
Set-TextCell $ws "E34" "  +0.63%  "
Set-TextCell $ws "D35" "50.53"
Set-TextCell $ws "E35" "  -1.03%  "
Set-TextCell $ws "D36" "0.0444"
Set-TextCell $ws "E36" "  +2.39%  "
Set-TextCell $ws "E37" "  -0.06%  "
Set-TextCell $ws "E38" "  +5.99%  "
Set-TextCell $ws "D39" "17.13"
Set-TextCell $ws "E39" "  +4.31%  "
Set-TextCell $ws "D40" "0.286"
Set-TextCell $ws "E40" "  +5.71%  "
Set-TextCell $ws "D41" "1.88"
Set-TextCell $ws "E41" "  +3.20%  "
Set-TextCell $ws "D42" "128.67"
Set-TextCell $ws "E42" "  +3.22%  "
Set-TextCell $ws "E43" "  +0.27%  "
Set-TextCell $ws "D44" "2.56"
Set-TextCell $ws "E44" "  +2.50%  "
Set-TextCell $ws "D45" "3.71"
Set-TextCell $ws "E45" "  +5.31%  "
Set-TextCell $ws "D46" "22.35"
Set-TextCell $ws "E46" "  +4.68%  "
Set-TextCell $ws "D47" "2.55"
Set-TextCell $ws "E47" "  +8.02%  "
Set-TextCell $ws "E48" "  +3.75%  "
Set-TextCell $ws "D49" "2.049.04"
Set-TextCell $ws "E49" "  +0.52%  "
Set-TextCell $ws "D50" "3.367.57"
Set-TextCell $ws "E50" "  +3.61%  "
Set-TextCell $ws "B51" "BEAM"
Set-TextCell $ws "C51" "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
Set-TextCell $ws "D51" "0.0319"
Set-TextCell $ws "E51" "  -0.62%  "
